# Output folder feature - extra context can be saved - repeated analyses are saved separately.
#
# 1) "Input" sheet: record the source filename in a new column D, widen the column.
# 2) "Calc" sheet: add two new "Montefehler Erfolgsrate" (Monte Carlo error success rate)
#    columns (BG/BH), tweak a handful of re-computed values, adjust a few column widths.
# 3) "Results" sheet: mirrors the recomputed age values, a couple of column widths change.
# 4) "Constants" sheet: R30_29 constant refined from 4.8E-05 to 5E-05.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Input sheet
# ---------------------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Range("D3").Value = "Filename: C:\Neptune\User\Neptune\Data\UTh\2015\0815\011_7184.dat"
$wsInput.Columns.Item(4).ColumnWidth = 65.87760416666667

# ---------------------------------------------------------------------------
# Calc sheet
# ---------------------------------------------------------------------------
$wsCalc = $wb.Worksheets.Item("Calc")

# New columns BG / BH: "Unkorr. Montefehler Erfolgsrate" / "Korr. Montefehler Erfolgsrate"
$wsCalc.Range("BG1").Value = "Unkorr. Montefehler Erfolgsrate"
$wsCalc.Range("BG1").Font.Bold = $true
$wsCalc.Range("BG1").HorizontalAlignment = -4108

$wsCalc.Range("BH1").Value = "Korr. Montefehler Erfolgsrate"
$wsCalc.Range("BH1").Font.Bold = $true
$wsCalc.Range("BH1").HorizontalAlignment = -4108

$wsCalc.Range("BG2").Value = "(%)"
$wsCalc.Range("BG2").Font.Bold = $true
$wsCalc.Range("BG2").HorizontalAlignment = -4108
$wsCalc.Range("BG2").Borders.Item(9).LineStyle = -4119

$wsCalc.Range("BH2").Value = "(%)"
$wsCalc.Range("BH2").Font.Bold = $true
$wsCalc.Range("BH2").HorizontalAlignment = -4108
$wsCalc.Range("BH2").Borders.Item(9).LineStyle = -4119

$wsCalc.Range("BG3").Value = 100
$wsCalc.Range("BH3").Value = 100

# Recomputed values (re-run of the Monte Carlo error propagation)
$wsCalc.Range("AP3").Value = 0.5328000000000001
$wsCalc.Range("AQ3").Value = 0.1886865953897859
$wsCalc.Range("AW3").Value = 0.5317
$wsCalc.Range("AX3").Value = 0.5340553074276952
$wsCalc.Range("AY3").Value = 0.1883322553611112
$wsCalc.Range("BC3").Value = 0.5688238147912374
$wsCalc.Range("BE3").Value = 267.0276537138476
$wsCalc.Range("BF3").Value = 0.1891665234820942

# Column width tweaks
$wsCalc.Columns.Item(49).ColumnWidth = 8.877604166666666
$wsCalc.Columns.Item(55).ColumnWidth = 18.877604166666668
$wsCalc.Columns.Item(59).ColumnWidth = 31.877604166666668
$wsCalc.Columns.Item(60).ColumnWidth = 29.877604166666668

# ---------------------------------------------------------------------------
# Results sheet
# ---------------------------------------------------------------------------
$wsResults = $wb.Worksheets.Item("Results")
$wsResults.Range("N3").Value = 0.5328000000000001
$wsResults.Range("P3").Value = 0.5317
$wsResults.Range("R3").Value = 0.5688238147912374

$wsResults.Columns.Item(16).ColumnWidth = 7.877604166666667
$wsResults.Columns.Item(18).ColumnWidth = 18.877604166666668

# ---------------------------------------------------------------------------
# Constants sheet
# ---------------------------------------------------------------------------
$wsConstants = $wb.Worksheets.Item("Constants")
$wsConstants.Range("B3").Value = 0.00005
